# Update column G (K = strikeouts) values for rows 2-37 in the active worksheet.
# The commit regenerates the save_data to use K instead of Strike# and
# recalculates the derived stats; here we directly write the new K (G column)
# values per row as captured from the regenerated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 1
    4  = 1
    5  = 3
    6  = 1
    7  = 3
    8  = 2
    9  = 2
    10 = 2
    11 = 1
    12 = 2
    13 = 2
    14 = 3
    15 = 1
    16 = 3
    17 = 4
    18 = 2
    19 = 1
    20 = 2
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 2
    27 = 1
    28 = 1
    29 = 2
    30 = 1
    31 = 1
    32 = 2
    33 = 2
    34 = 2
    35 = 2
    36 = 1
    37 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
